$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.956.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "2.234.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.19%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'304.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.34%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'95.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.52%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.569"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.02%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D9").Value = "'0.521"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.95%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'34.81"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.59%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0807"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.80%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "  -5.87%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.104"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.88%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "2.573.56"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.34%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "2.274.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.28%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.822"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.59%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'13.57"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.22%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "43.827.12"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "0.0₃0958"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.23%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'12.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -9.72%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'6.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.86%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'64.82"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.49%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'235.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -6.39%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.17%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "  -7.40%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'9.91"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.48%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'37.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.02%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D30").Value = "'5.91"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.80%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'19.87"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Value = "'154.69"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.82%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.0801"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.81%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'3.20"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.99%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'2.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.73%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "  -1.26%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "  -6.38%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "  -12.18%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'15.17"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.79%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = "NEARProtocol"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'3.34"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.74%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "RenderToken"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'3.80"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -9.76%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "  -5.75%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "1.742.11"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.65%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'85.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.85%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "  -5.45%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "  -5.22%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E49").Value = "  -8.11%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'8.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.82%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'54.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.24%  "
$ws.Range("E51").Style = "Normal"
